$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.801.95'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '2.341.47'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '539.18'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.04'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('E7').Value = '  +0.71%  '
$ws.Range('E8').Value = '  +6.17%  '
$ws.Range('E9').Value = '  +0.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.49'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.08%  '
$ws.Range('E11').Value = '  -1.69%  '
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('E13').Value = '  +0.86%  '
$ws.Range('D14').Value = '2.758.98'
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('D15').Value = '57.746.82'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('D17').Value = '2.345.36'
$ws.Range('E17').Value = '  -0.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.68'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.30'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '328.36'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.69'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '62.73'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.54%  '
$ws.Range('E24').Value = '  -3.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.997'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.31'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.44%  '
$ws.Range('E27').Value = '  -5.99%  '
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '169.87'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.41%  '
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('B32').Value = 'SuiNetwork'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.02'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.31'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.37%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('E35').Value = '  +0.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.18'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.44%  '
$ws.Range('E37').Value = '  -2.04%  '
$ws.Range('E38').Value = '  -0.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '39.09'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '141.44'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -4.91%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.376'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.48%  '
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '287.42'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0946'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0510'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.12'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.92%  '
$ws.Range('E47').Value = '  +1.08%  '
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.379'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.08'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.66%  '
$ws.Range('E51').Value = '  +0.97%  '
